# Generate Report for Handoff
# Adds a new "handback" row (row 3) to the Overview, zh-cn and de-de sheets
# describing the file 20c5395d-0b78-4783-9f40-8d8f07e67efeoo....md

$wb = $excel.ActiveWorkbook

$newFileName    = '20c5395d-0b78-4783-9f40-8d8f07e67efeooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$newFilePath    = 'e2e\20c5395d-0b78-4783-9f40-8d8f07e67efeooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$readyStatus    = 'Ready for handoff'
$handoffDate    = '2016-08-21 20:40:20'
$handoffDateZh  = '2016-08-21 20:40:15'
$ext            = '.md'
$blank          = ''
$e2e            = 'e2e'
$ht             = 'ht'
$falseStr       = 'False'
$trueStr        = 'True'
$epoch          = '0001-01-01 00:00:00'
$zhXlf          = '20c5395d-0b78-4783-9f40-8d8f07e67efeoooooooooooooooooooooooooooooooooooooooo.5da184f345abf8274b2d8f9f45b07ed825a9200c.zh-cn.xlf'
$deXlf          = '20c5395d-0b78-4783-9f40-8d8f07e67efeoooooooooooooooooooooooooooooooooooooooo.5da184f345abf8274b2d8f9f45b07ed825a9200c.de-de.xlf'

$hyperlinkUrl = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a11ec39c8be0bf076a457de74d8ad876baf8ba0a/e2e/20c5395d-0b78-4783-9f40-8d8f07e67efeooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1 / table3): columns A-G
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newFileName
$wsOverview.Range("B3").Value = $newFilePath
$wsOverview.Range("C3").Value = $ext
$wsOverview.Range("D3").Value = $blank
$wsOverview.Range("E3").Value = $readyStatus
$wsOverview.Range("F3").Value = $readyStatus
$wsOverview.Range("G3").Value = $handoffDate
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkUrl, "", "", $newFilePath) | Out-Null
$wsOverview.Range("B3").Font.Underline = 2
$wsOverview.Range("B3").Font.Color = 6495237

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2 / table1): columns A-P
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = $newFileName
$wsZh.Range("B3").Value = $ext
$wsZh.Range("C3").Value = $readyStatus
$wsZh.Range("D3").Value = $e2e
$wsZh.Range("E3").Value = $ht
$wsZh.Range("F3").Value = $falseStr
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = $handoffDateZh
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I3").Value = $blank
$wsZh.Range("J3").Value = $blank
$wsZh.Range("K3").Value = $epoch
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L3").Value = $blank
$wsZh.Range("M3").Value = $trueStr
$wsZh.Range("N3").Value = $blank
$wsZh.Range("O3").Value = $falseStr
$wsZh.Range("P3").Value = $blank

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $hyperlinkUrl, "", "", $newFileName) | Out-Null
$wsZh.Range("A3").Font.Underline = 2
$wsZh.Range("A3").Font.Color = 6495237

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3 / table2): columns A-P
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = $newFileName
$wsDe.Range("B3").Value = $ext
$wsDe.Range("C3").Value = $readyStatus
$wsDe.Range("D3").Value = $e2e
$wsDe.Range("E3").Value = $ht
$wsDe.Range("F3").Value = $falseStr
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = $handoffDate
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I3").Value = $blank
$wsDe.Range("J3").Value = $blank
$wsDe.Range("K3").Value = $epoch
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L3").Value = $blank
$wsDe.Range("M3").Value = $trueStr
$wsDe.Range("N3").Value = $blank
$wsDe.Range("O3").Value = $falseStr
$wsDe.Range("P3").Value = $blank

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $hyperlinkUrl, "", "", $newFileName) | Out-Null
$wsDe.Range("A3").Font.Underline = 2
$wsDe.Range("A3").Font.Color = 6495237

# ---------------------------------------------------------------------------
# Column width adjustments (E/F on Overview, C on zh-cn/de-de widen slightly)
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 17.08
$wsOverview.Columns.Item(6).ColumnWidth = 17.08
$wsZh.Columns.Item(3).ColumnWidth = 17.08
$wsDe.Columns.Item(3).ColumnWidth = 17.08
